# Add five "blank for whiteboard" slides after the existing title slide.
# Slide 2 is built from scratch on the Blank layout; slides 3-6 are exact
# duplicates of slide 2 (matching how the author repeatedly duplicated the
# slide in PowerPoint).

$p = $ppt.ActivePresentation

# ppLayoutBlank = 12 -> slideLayout12.xml ("Blank") in this deck's master.
$s = $p.Slides.Add(2, 12)
$s.SlideShowTransition.Speed = 2  # ppTransitionSpeedMedium -> spd="med"

$tb = $s.Shapes.AddTextbox(1, 0, 0, 100, 100)
$tb.Name = "TextBox 5"

$tf = $tb.TextFrame
$tr = $tf.TextRange
$tr.Text = "This slide left blank for whiteboard"

# Paragraph is centered.
$tr.ParagraphFormat.Alignment = 2

# Only the first three characters ("Thi") carry explicit run formatting;
# the remainder inherits the (identical-looking) default text style.
$firstRun = $tr.Characters(1, 3)
$firstRun.Font.Size = 30
$firstRun.Font.Bold = $true
$firstRun.Font.NameAscii = "Helvetica Neue"
$firstRun.Font.NameFarEast = "Helvetica Neue"
$firstRun.Font.NameComplexScript = "Helvetica Neue"
$firstRun.Font.Color.RGB = 0

# Autofit the box to the text, center it vertically, 4pt insets on all sides.
$tf.WordWrap = $true
$tf.AutoSize = 1
$tf.VerticalAnchor = 3
$tf.MarginLeft = 4
$tf.MarginRight = 4
$tf.MarginTop = 4
$tf.MarginBottom = 4

# No fill, hairline (1pt) invisible outline.
$tb.Fill.Visible = $false
$tb.Line.Visible = $false
$tb.Line.Weight = 1

# Final position/size in points (source EMU / 12700).
$tb.Left = 491.11110236220475
$tb.Top = 1002.2296062992126
$tb.Width = 1031.1111023622047
$tb.Height = 44.429686

# Slides 3-6 are plain duplicates of slide 2, added one after another so
# the resulting slide order stays 1,2,3,4,5,6.
for ($i = 0; $i -lt 4; $i++) {
    $last = $p.Slides.Item($p.Slides.Count)
    $last.Duplicate() | Out-Null
}

Write-Host "Slide count: $($p.Slides.Count)"
